# Auto-generated edit script
# Updates cryptos list price (D) and volume change (E) columns per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.888.31"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "1.584.35"
$ws.Range("E3").Value = "  -1.97%  "
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").Value = "210.19"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").Value = "0.477"
$ws.Range("E7").Value = "  -2.05%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("D10").Value = "18.09"
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("D11").Value = "0.0790"
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").Value = "1.805.45"
$ws.Range("E12").Value = "  -1.89%  "
$ws.Range("D13").Value = "1.584.07"
$ws.Range("E13").Value = "  -2.00%  "
$ws.Range("E14").Value = "  -2.18%  "
$ws.Range("E15").Value = "  -2.21%  "
$ws.Range("D16").Value = "25.886.52"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("E17").Value = "  -1.00%  "
$ws.Range("D18").Value = "59.98"
$ws.Range("E18").Value = "  -2.42%  "
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").Value = "193.03"
$ws.Range("E20").Value = "  +1.01%  "
$ws.Range("D21").Value = "4.19"
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("E22").Value = "  -0.59%  "
$ws.Range("E23").Value = "  -0.97%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").Value = "141.44"
$ws.Range("E25").Value = "  -1.14%  "
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").Value = "15.08"
$ws.Range("E28").Value = "  -0.35%  "
$ws.Range("E29").Value = "  -2.41%  "
$ws.Range("E30").Value = "  -4.98%  "
$ws.Range("D31").Value = "0.0471"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("E33").Value = "  -1.59%  "
$ws.Range("E34").Value = "  +1.06%  "
$ws.Range("E35").Value = "  -2.15%  "
$ws.Range("D36").Value = "1.096.13"
$ws.Range("E36").Value = "  -2.60%  "
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("E38").Value = "  -1.85%  "
$ws.Range("E39").Value = "  -0.61%  "
$ws.Range("E40").Value = "  -2.64%  "
$ws.Range("E41").Value = "  -4.75%  "
$ws.Range("D42").Value = "0.797"
$ws.Range("E42").Value = "  +5.23%  "
$ws.Range("D43").Value = "93.45"
$ws.Range("E43").Value = "  -3.86%  "
$ws.Range("D44").Value = "5.13"
$ws.Range("E44").Value = "  +1.26%  "
$ws.Range("D45").Value = "1.719.11"
$ws.Range("E45").Value = "  -1.83%  "
$ws.Range("E46").Value = "  -1.80%  "
$ws.Range("E47").Value = "  +2.02%  "
$ws.Range("D48").Value = "53.14"
$ws.Range("E48").Value = "  -1.03%  "
$ws.Range("D49").Value = "0.0508"
$ws.Range("E49").Value = "  -1.23%  "
$ws.Range("E50").Value = "  -0.97%  "
$ws.Range("E51").Value = "  -0.21%  "
